$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a value to be stored as text (shared string) even when it
# looks numeric (e.g. RUT numbers), matching how the source data was typed.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 3 becomes what used to be row 2's person (Maximiliano Joaquin Almonacid Perez)
Set-TextValue "B3" "21494146"
Set-TextValue "C3" "5"
Set-TextValue "D3" "MAXIMILIANO JOAQUIN"
Set-TextValue "E3" "ALMONACID PÉREZ"

# Row 2 becomes a new person (Flavio Alexander Jara Labrin)
Set-TextValue "B2" "21075353"
Set-TextValue "C2" "2"
Set-TextValue "D2" "FLAVIO ALEXANDER"
Set-TextValue "E2" "JARA LABRIN"

# New row 4 is added for another new person (Mathias Eduardo Deumacan Pulgar)
Set-TextValue "A4" "31-12-2024"
Set-TextValue "B4" "21223313"
Set-TextValue "C4" "8"
Set-TextValue "D4" "MATHIAS EDUARDO"
Set-TextValue "E4" "DEUMACÁN PULGAR"
Set-TextValue "F4" "1"
Set-TextValue "G4" "FÍSICA MECANICA / 3"
